$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 269-272 (columns H and I) with corrected AgTests/AgPosit figures
$ws.Range("H269").Value = 9420
$ws.Range("I269").Value = 409

$ws.Range("H270").Value = 2649
$ws.Range("I270").Value = 176

$ws.Range("H271").Value = 50921
$ws.Range("I271").Value = 1816

$ws.Range("H272").Value = 30042
$ws.Range("I272").Value = 1668

# Give the two new date cells the same style (date number format) as the rest of column A
$ws.Range("A272").Copy()
$ws.Range("A273:A274").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 273 (2020-12-02 daily stats)
$ws.Range("A273").Value = 44167
$ws.Range("B273").Value = 111208
$ws.Range("C273").Value = 74273
$ws.Range("D273").Value = 36005
$ws.Range("E273").Value = 10264
$ws.Range("F273").Value = 1982
$ws.Range("G273").Value = 930
$ws.Range("H273").Value = 24928
$ws.Range("I273").Value = 1244

# Add new row 274 (2020-12-03 daily stats)
$ws.Range("A274").Value = 44168
$ws.Range("B274").Value = 113392
$ws.Range("C274").Value = 77142
$ws.Range("D274").Value = 35293
$ws.Range("E274").Value = 11221
$ws.Range("F274").Value = 2184
$ws.Range("G274").Value = 957
$ws.Range("H274").Value = 15794
$ws.Range("I274").Value = 628

# Update the view: scroll so row 238 is at the top, and select G271
$ws.Application.ActiveWindow.ScrollRow = 238
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G271").Select()
